$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and 1h volume change (column E) values
$ws.Range("D2").Formula = "'26.222.12"
$ws.Range("D3").Formula = "'1.588.63"
$ws.Range("E3").Value = "  +0.79%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Formula = "'212.20"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").Formula = "'0.501"
$ws.Range("E6").Value = "  +0.93%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").Formula = "'19.35"
$ws.Range("E10").Value = "  -0.87%  "
$ws.Range("D11").Formula = "'0.0849"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Formula = "'1.811.24"
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("D13").Formula = "'1.601.60"
$ws.Range("E13").Value = "  +1.32%  "
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Formula = "'26.234.22"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("D18").Formula = "'0.0₃0726"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("E19").Value = "  +1.87%  "
$ws.Range("D20").Formula = "'213.37"
$ws.Range("E20").Value = "  +3.10%  "
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").Formula = "'8.99"
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("E24").Value = "  -2.45%  "
$ws.Range("D25").Formula = "'143.97"
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("D27").Formula = "'7.05"
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("E30").Value = "  -1.88%  "
$ws.Range("D31").Formula = "'1.15"
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("D32").Formula = "'3.19"
$ws.Range("D33").Formula = "'2.94"
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("D34").Formula = "'1.332.24"
$ws.Range("E34").Value = "  +4.12%  "
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("D37").Formula = "'0.589"
$ws.Range("E37").Value = "  -3.59%  "
$ws.Range("D38").Formula = "'0.0166"
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("D40").Formula = "'5.73"
$ws.Range("E40").Value = "  +2.91%  "
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").Formula = "'1.01"
$ws.Range("E42").Value = "  -7.45%  "
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").Formula = "'0.765"
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("D45").Formula = "'61.81"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("D46").Formula = "'1.723.14"
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("D47").Formula = "'85.45"
$ws.Range("E47").Value = "  -4.21%  "
$ws.Range("D48").Formula = "'1.48"
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("D50").Formula = "'0.0972"
$ws.Range("E50").Value = "  -3.01%  "
$ws.Range("E51").Value = "  -0.23%  "
